$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Sheet1 (quality_comparison)
$ws1.Range("C2").Value = "approach"

# Sheet2 (computational_comparison)
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Borders for C1/D1 on sheet1
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1
